$d = $word.ActiveDocument

# 1) "...generated with a web app..." -> "...generated by a web app..."
$r = $d.Content
$r.Find.Execute("generated with a web app", $true, $false, $false, $false, $false, $true, 1, $false, "generated by a web app", 2) | Out-Null

# 2) "...with the following parameters:" -> "...with the following input parameters:"
$r = $d.Content
$r.Find.Execute("with the following parameters:", $true, $false, $false, $false, $false, $true, 1, $false, "with the following input parameters:", 2) | Out-Null

# 3) "...which designed and populate two circuit boards, used different ICs on each of them. So this cannot be correct. "
#    -> "...which designed and populated two circuit boards, used different ICs on each of them. So at least one of
#       the boards cannot work. They also told us that when they tried to figure out the error, they probably
#       destroyed some parts. "
$r = $d.Content
$r.Find.Execute( `
    "designed and populate two circuit boards, used different ICs on each of them. So this cannot be correct. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "designed and populated two circuit boards, used different ICs on each of them. So at least one of the boards cannot work. They also told us that when they tried to figure out the error, they probably destroyed some parts. ", `
    2) | Out-Null

# 4) "Other big error sources" -> "Big error sources"
$r = $d.Content
$r.Find.Execute("Other big error sources", $true, $false, $false, $false, $false, $true, 1, $false, "Big error sources", 2) | Out-Null

# 5) Split the paragraph that ends with "...on the plug board." so the trailing
#    _GoBack bookmark sits in its own, new, empty paragraph.
$r = $d.Content
$r.Find.Execute("on the plug board.", $true, $false, $false, $false, $false, $true, 1, $false, "on the plug board.^p", 2) | Out-Null

# 6) Remove the "Improvement suggestions for the last board" paragraph together
#    with the trailing empty paragraph that used to follow it.
$r = $d.Content
$found = $r.Find.Execute("Improvement suggestions for the last board")
if ($found) {
    $start = $r.Start
    $end = $d.Content.End
    $d.Range($start, $end).Delete() | Out-Null
}
